$p = $ppt.ActivePresentation

# Slide 15: title "阅读经文：路加福音9: 52-56；10: 1-42" -> "阅读经文：路加福音"
# Remove the scripture-reference runs that follow "路加福音" (chars 10-26).
$s15 = $p.Slides.Item(15)
$titleShape15 = $s15.Shapes.Item(1)
$tr15 = $titleShape15.TextFrame.TextRange
$tr15.Characters(10, 17).Text = ""

# Slide 16: title "证道经文：路加福音10: 38-42" -> "证道经文：路加福音"
# Remove the scripture-reference runs that follow "路加福音" (chars 10-18).
$s16 = $p.Slides.Item(16)
$titleShape16 = $s16.Shapes.Item(1)
$tr16 = $titleShape16.TextFrame.TextRange
$tr16.Characters(10, 9).Text = ""

# Slide 16: body placeholder - two leading-space runs become empty runs.
$bodyShape16 = $s16.Shapes.Item(2)
$trBody16 = $bodyShape16.TextFrame.TextRange
$trBody16.Characters(1, 1).Text = ""
# After removing the first space, the second space (originally char 9) is now char 8.
$trBody16.Characters(8, 1).Text = ""
